$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume(1h) (E) columns for rows 2-51
# so numeric-looking strings are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.712.96"
$ws.Range("E2").Value = "  +11.28%  "
$ws.Range("D3").Value = "1.685.14"
$ws.Range("E3").Value = "  +6.24%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "306.16"
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("D6").Value = "0.9958"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "0.3685"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D8").Value = "49.53"
$ws.Range("E8").Value = "  +20.37%  "
$ws.Range("D9").Value = "0.3425"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "1.164"
$ws.Range("E10").Value = "  +4.79%  "
$ws.Range("D11").Value = "0.07237"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("D12").Value = "0.9977"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "6.098"
$ws.Range("E13").Value = "  +5.48%  "
$ws.Range("D14").Value = "20.17"
$ws.Range("E14").Value = "  +4.64%  "
$ws.Range("D15").Value = "6.698"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "1.682.13"
$ws.Range("E16").Value = "  +5.99%  "
$ws.Range("D17").Value = "0.00001103"
$ws.Range("E17").Value = "  +4.27%  "
$ws.Range("D18").Value = "0.9953"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "0.06664"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "80.94"
$ws.Range("E20").Value = "  +6.92%  "
$ws.Range("D21").Value = "16.37"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("D22").Value = "6.077"
$ws.Range("E22").Value = "  +3.23%  "
$ws.Range("E23").Value = "  +4.93%  "
$ws.Range("D24").Value = "24.633.57"
$ws.Range("E24").Value = "  +11.14%  "
$ws.Range("D25").Value = "2.412"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "3.346"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.663"
$ws.Range("E27").Value = "  +7.11%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "152.41"
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.44"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.868.40"
$ws.Range("E30").Value = "  +6.35%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "127.59"
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "6.257"
$ws.Range("E32").Value = "  +8.32%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "4.030"
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.9793"
$ws.Range("E34").Value = "  +7.63%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.08409"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "1.686"
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "12.34"
$ws.Range("E37").Value = "  +6.18%  "
$ws.Range("D38").Value = "0.06347"
$ws.Range("E38").Value = "  +6.27%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.302"
$ws.Range("E39").Value = "  +4.56%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02308"
$ws.Range("E40").Value = "  +6.48%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "8.639"
$ws.Range("E41").Value = "  +4.63%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.244"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.2084"
$ws.Range("E43").Value = "  +6.04%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.6090"
$ws.Range("E44").Value = "  +6.17%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "0.9954"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "13.09"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.750"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.5878"
$ws.Range("E48").Value = "  +6.18%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "125.48"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.998"
$ws.Range("E50").Value = "  +3.70%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.07236"
$ws.Range("E51").Value = "  +7.94%  "
